# Applies:
#  1. Adds <w:spacing w:line="360" w:lineRule="auto"/> to the pPr of the
#     "It is recommended to install motion sensors..." paragraph.
#  2. Splits "${IC}" into three runs "${" / "M" / "IC}" (-> ${MIC}).
#  3. Splits "${PB}" into three runs "${" / "M" / "PB}" (-> ${MPB}).

$d = $word.ActiveDocument

# --- 1. Paragraph line spacing -------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*It is recommended to install motion sensor*") {
        $p.Format.LineSpacingRule = 5    # wdLineSpaceMultiple
        $p.Format.LineSpacing = 18       # 1.5 lines (360 twips / 240 * 12pt)
        break
    }
}

# --- helper: split "${TAG}" into "${" / "M" / "TAG}" inside a table cell -
function Split-RebatePlaceholder($table, $row, $col, $tag) {
    $cell = $table.Cell($row, $col)
    $cellStart = $cell.Range.Start
    $full = "`${$tag}"
    $len = $full.Length

    # Isolate the "${" prefix from the "TAG}" suffix by toggling Bold off/on
    # on the suffix range - this forces Word to split the run without
    # altering any text, so the original run (and its rsid) survives on the
    # untouched "${" prefix.
    $suffix = $d.Range($cellStart + 2, $cellStart + $len)
    $suffix.Bold = 1
    $suffix2 = $d.Range($cellStart + 2, $cellStart + $len)
    $suffix2.Bold = 0

    # Retext the suffix run "TAG}" -> "MTAG}", protecting the prefix run
    # with a temporary Bold toggle so the two don't get recombined (which
    # would also wipe the prefix run's rsid).
    $prefixGuard1 = $d.Range($cellStart, $cellStart + 2)
    $prefixGuard1.Bold = 1

    $suffixRange = $d.Range($cellStart + 2, $cellStart + $len)
    $suffixRange.Text = "M$tag}"

    $prefixGuard2 = $d.Range($cellStart, $cellStart + 2)
    $prefixGuard2.Bold = 0

    # Split the now "M" + "TAG}" run into its own two runs by toggling Bold
    # on just the "M" character. Guard the prefix again so it stays put.
    $prefixGuard3 = $d.Range($cellStart, $cellStart + 2)
    $prefixGuard3.Bold = 1

    $mChar = $d.Range($cellStart + 2, $cellStart + 3)
    $mChar.Bold = 1
    $mChar2 = $d.Range($cellStart + 2, $cellStart + 3)
    $mChar2.Bold = 0

    $prefixGuard4 = $d.Range($cellStart, $cellStart + 2)
    $prefixGuard4.Bold = 0
}

$table = $d.Tables.Item(1)
Split-RebatePlaceholder $table 2 2 "IC"
Split-RebatePlaceholder $table 3 2 "PB"
